$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" column (E) for the data rows (16-51) is reversed in
# order: what used to be the ascending chronological list (1810 .. 2109)
# is now listed descending (2109 .. 1810), while every other column,
# style and row stays exactly where it was.
$periods = @(
    "2109","2108","2107","2106","2105","2104","2103","2102","2101",
    "2012","2011","2010","2009","2008","2007","2006","2005","2004","2003","2002","2001",
    "1912","1911","1910","1909","1908","1907","1906","1905","1904","1903","1902","1901",
    "1812","1811","1810"
)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("E$row").Value = $periods[$i]
}
